$d = $word.ActiveDocument

function Find-ParagraphIndex($pattern) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -match $pattern) {
            return $i
        }
    }
    return -1
}

function Set-ParagraphXml($index, $innerXml) {
    $p = $d.Paragraphs.Item($index)
    $rng = $p.Range
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xml) | Out-Null
}

# 1. Insert a new checklist paragraph right after "__ Open experiment"
$idx = Find-ParagraphIndex("^__ Open experiment")
$p = $d.Paragraphs.Item($idx)
$p.Range.InsertParagraphAfter()
$newP = $d.Paragraphs.Item($idx + 1)
$newP.Range.Text = "__ Look up subject number (for return visits) or assign new number (for first visits)"

# 2. cd ~/Desktop/x121 -> cd ~/Desktop/x122 (rewritten via InsertXML so the
#    leading <w:tab/> run child survives -- a plain Find/Replace would fold
#    the tab into the text as a literal \t character)
$idx = Find-ParagraphIndex("cd ~/Desktop/x121")
$innerXml = '<w:p><w:r><w:tab/><w:t>cd ~/Desktop/x122</w:t></w:r></w:p>'
Set-ParagraphXml $idx $innerXml

# 3. Rewrite the ./exp.py line into its new multi-run form
$idx = Find-ParagraphIndex("^\t?\./exp\.py")
$innerXml = '<w:p><w:r><w:tab/><w:t xml:space="preserve">./exp.py </w:t></w:r><w:r><w:t>#</w:t></w:r><w:r><w:t xml:space="preserve">         </w:t></w:r><w:r><w:t xml:space="preserve">  </w:t></w:r><w:r><w:t xml:space="preserve">                      (dot slash exp dot py</w:t></w:r><w:r><w:t xml:space="preserve"> space &lt;</w:t></w:r><w:r><w:t xml:space="preserve">subject  </w:t></w:r><w:r><w:t>#&gt;</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p>'
Set-ParagraphXml $idx $innerXml

# 4. Rewrite the "$25" paragraph text and move the _GoBack bookmark into it
$idx = Find-ParagraphIndex("^You will also have the opportunity")
$innerXml = '<w:p><w:r><w:t>You will also have the opportunity to earn $25 based on your participation. During the experime</w:t></w:r><w:r><w:t>nt, you will accumulate points. The highest score for each of the three sessions will recieve</w:t></w:r><w:r><w:t xml:space="preserve"> will each receive $25 after the experiment is completed in a few weeks.</w:t></w:r><w:r><w:t xml:space="preserve"> (Maximum $25 per participant).</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
Set-ParagraphXml $idx $innerXml

# 5. Add a lastRenderedPageBreak inside the "Please make sure..." paragraph
$idx = Find-ParagraphIndex("^Please make sure")
$innerXml = '<w:p><w:r><w:t xml:space="preserve">Please make sure to read all the instructions CAREFULLY before each block. Different blocks of the experiment may have different instructions: Sometimes you </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>will be given feedback</w:t></w:r><w:r><w:t xml:space="preserve"> after an incorrect response, other times no feedback will be given.</w:t></w:r></w:p>'
Set-ParagraphXml $idx $innerXml

# 6. Remove the "Write down participant score..." paragraph entirely
$idx = Find-ParagraphIndex("^__ Write down participant score")
$d.Paragraphs.Item($idx).Range.Delete() | Out-Null

# 7. Remove the "Have participant take questionnaire" paragraph entirely
$idx = Find-ParagraphIndex("^__ Have participant take questionnaire")
$d.Paragraphs.Item($idx).Range.Delete() | Out-Null

# 8. Remove the "Write subject ID number on questionnaire" paragraph (it used to
#    carry the _GoBack bookmark, which now lives in the "$25" paragraph instead)
$idx = Find-ParagraphIndex("^__ Write subject ID number on questionnaire")
$d.Paragraphs.Item($idx).Range.Delete() | Out-Null

Write-Output "Edits applied. Final paragraph count: $($d.Paragraphs.Count)"
